# Generate the YCbCr color-space GLCM data (135-degree, R channel) on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(0.11619532427520814, 0.75928194804744686, 0.5316018625978477, 0.94329419823224059),
  @(0.25316684870370409, 0.63241344154999701, 0.34892987116351504, 0.8798100093697212),
  @(0.060866522902294967, 0.73360624060671598, 0.73779191839507163, 0.9698900301951433),
  @(0.13617758289594026, 0.81187198316736453, 0.38259127928498105, 0.93270931154342951)
)

for ($r = 0; $r -lt 4; $r++) {
  for ($c = 0; $c -lt 4; $c++) {
    $ws.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
  }
}

# Columns A:D get a fixed custom width, matching the authored sheet.
$ws.Columns("A:D").ColumnWidth = 12.7109375

# Keep the sheet marked as the selected tab (same as the source file).
$ws.Tab.Selected = $true

# Pick up the text / date-time number-format styles used alongside this
# data (registered in the workbook's style table, then cleared again so
# the worksheet's used range stays A1:D4 with plain numeric cells).
$ws.Range("F1").NumberFormat = "@"
$ws.Range("G1").NumberFormat = "m/d/yy h:mm"
$ws.Range("F1:G1").Clear()

# Make sure the workbook fully recalculates the next time it is opened.
$excel.CalculateFullRebuild()
